$wb = $excel.ActiveWorkbook

# Crime data update for 2022-12-13 (new records added across citywide,
# "By Neighborhood" summary and per-neighborhood detail sheets, column I = 2022).

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("I2").Value = 6969
$ws.Range("I3").Value = 7213
$ws.Range("I4").Value = 1657
$ws.Range("I5").Value = 678
$ws.Range("I6").Value = 8514
$ws.Range("I7").Value = 25031

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("I2").Value = 199
$ws.Range("I6").Value = 182
$ws.Range("I7").Value = 785
$ws.Range("I8").Value = 1486
$ws.Range("I9").Value = 132
$ws.Range("I14").Value = 140
$ws.Range("I15").Value = 292
$ws.Range("I16").Value = 74
$ws.Range("I19").Value = 702
$ws.Range("I20").Value = 622
$ws.Range("I23").Value = 246
$ws.Range("I31").Value = 253
$ws.Range("I33").Value = 1101
$ws.Range("I36").Value = 341
$ws.Range("I37").Value = 773
$ws.Range("I42").Value = 942
$ws.Range("I51").Value = 293
$ws.Range("I52").Value = 568
$ws.Range("I55").Value = 289
$ws.Range("I63").Value = 77
$ws.Range("I65").Value = 583
$ws.Range("I67").Value = 947
$ws.Range("I72").Value = 99
$ws.Range("I76").Value = 356
$ws.Range("I79").Value = 717
$ws.Range("I83").Value = 539
$ws.Range("I85").Value = 1114
$ws.Range("I88").Value = 231
$ws.Range("I90").Value = 326
$ws.Range("I91").Value = 261
$ws.Range("I93").Value = 144
$ws.Range("I95").Value = 386
$ws.Range("I96").Value = 294
$ws.Range("I97").Value = 228
$ws.Range("I99").Value = 437
$ws.Range("I101").Value = 25031

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("I2").Value = 318
$ws.Range("I3").Value = 419
$ws.Range("I7").Value = 1114

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("I6").Value = 187
$ws.Range("I7").Value = 568

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("I2").Value = 439
$ws.Range("I7").Value = 1486

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("I2").Value = 258
$ws.Range("I3").Value = 238
$ws.Range("I6").Value = 213
$ws.Range("I7").Value = 785

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range("I2").Value = 74
$ws.Range("I3").Value = 68

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range("I3").Value = 67
$ws.Range("I7").Value = 294

$ws = $wb.Worksheets.Item('Bridgeport')
$ws.Range("I3").Value = 35
$ws.Range("I7").Value = 140

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("I2").Value = 229
$ws.Range("I3").Value = 251
$ws.Range("I7").Value = 773

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range("I2").Value = 124
$ws.Range("I7").Value = 437

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("I2").Value = 226
$ws.Range("I4").Value = 57
$ws.Range("I7").Value = 947

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range("I2").Value = 76
$ws.Range("I7").Value = 253

$ws = $wb.Worksheets.Item('New City')
$ws.Range("I3").Value = 176
$ws.Range("I7").Value = 583

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("I2").Value = 179
$ws.Range("I3").Value = 196
$ws.Range("I5").Value = 22
$ws.Range("I6").Value = 120
$ws.Range("I7").Value = 539

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("I4").Value = 18
$ws.Range("I6").Value = 78
$ws.Range("I7").Value = 386

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("I2").Value = 248
$ws.Range("I6").Value = 354
$ws.Range("I7").Value = 1101

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("I2").Value = 226
$ws.Range("I7").Value = 702

$ws = $wb.Worksheets.Item('River North')
$ws.Range("I2").Value = 73
$ws.Range("I7").Value = 356

$ws = $wb.Worksheets.Item('Ashburn')
$ws.Range("I3").Value = 43
$ws.Range("I7").Value = 182

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("I2").Value = 211
$ws.Range("I3").Value = 271
$ws.Range("I7").Value = 942

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range("I6").Value = 90
$ws.Range("I7").Value = 289

$ws = $wb.Worksheets.Item('Douglas')
$ws.Range("I5").Value = 8
$ws.Range("I7").Value = 246

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range("I3").Value = 95
$ws.Range("I7").Value = 261

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("I2").Value = 209
$ws.Range("I3").Value = 236
$ws.Range("I7").Value = 717

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("I2").Value = 176
$ws.Range("I7").Value = 622

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range("I3").Value = 115
$ws.Range("I7").Value = 341

$ws = $wb.Worksheets.Item('West Lawn')
$ws.Range("I2").Value = 40
$ws.Range("I3").Value = 36
$ws.Range("I6").Value = 60
$ws.Range("I7").Value = 144

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range("I2").Value = 85
$ws.Range("I7").Value = 292

$ws = $wb.Worksheets.Item('Avalon Park')
$ws.Range("I3").Value = 43
$ws.Range("I7").Value = 132

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range("I3").Value = 62
$ws.Range("I7").Value = 199

$ws = $wb.Worksheets.Item('West Town')
$ws.Range("I6").Value = 149
$ws.Range("I7").Value = 228

$ws = $wb.Worksheets.Item('United Center')
$ws.Range("I6").Value = 74
$ws.Range("I7").Value = 231

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range("I3").Value = 84
$ws.Range("I7").Value = 326

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("I6").Value = 119
$ws.Range("I7").Value = 293

$ws = $wb.Worksheets.Item('Old Town')
$ws.Range("I3").Value = 20
$ws.Range("I7").Value = 99

$ws = $wb.Worksheets.Item('Bucktown')
$ws.Range("I6").Value = 51
$ws.Range("I7").Value = 74
